$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- Step 1: clear old data range (keep headers/years) ---
$ws.Range("B2:M6").ClearContents()

# --- Step 2: write new header row + data, columns reordered ---
$ws.Cells.Item(1, 2).Value = "Unclassified"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(6, 2).Value = 0.04644247079960034

$ws.Cells.Item(1, 3).Value = "Trade, transportation, and utilities"
$ws.Cells.Item(2, 3).Value = 37.89317763408589
$ws.Cells.Item(3, 3).Value = 37.83979048224639
$ws.Cells.Item(4, 3).Value = 36.95590329548892
$ws.Cells.Item(5, 3).Value = 36.28191637022958
$ws.Cells.Item(6, 3).Value = 37.57198589864893

$ws.Cells.Item(1, 4).Value = "Public administration"
$ws.Cells.Item(2, 4).Value = 2.875315847148386
$ws.Cells.Item(3, 4).Value = 3.205996568304879
$ws.Cells.Item(4, 4).Value = 3.470839234951826
$ws.Cells.Item(5, 4).Value = 3.327914353277595
$ws.Cells.Item(6, 4).Value = 3.139513265636859

$ws.Cells.Item(1, 5).Value = "Professional and business services"
$ws.Cells.Item(2, 5).Value = 5.17556853581944
$ws.Cells.Item(3, 5).Value = 4.741262520747966
$ws.Cells.Item(4, 5).Value = 4.352773806197757
$ws.Cells.Item(5, 5).Value = 4.991871466858145
$ws.Cells.Item(6, 5).Value = 5.461638494043399

$ws.Cells.Item(1, 6).Value = "Other services"
$ws.Cells.Item(2, 6).Value = 1.568354098444573
$ws.Cells.Item(3, 6).Value = 1.715885487261764
$ws.Cells.Item(4, 6).Value = 1.735419631273921
$ws.Cells.Item(5, 6).Value = 2.075164958444865
$ws.Cells.Item(6, 6).Value = 2.136355193248264

$ws.Cells.Item(1, 7).Value = "Natural resources and mining"
$ws.Cells.Item(2, 7).Value = 1.062995555612427
$ws.Cells.Item(3, 7).Value = 1.164995979635621
$ws.Cells.Item(4, 7).Value = 1.25177809468938
$ws.Cells.Item(5, 7).Value = 1.262312324952634
$ws.Cells.Item(6, 7).Value = 1.300390117629378

$ws.Cells.Item(1, 8).Value = "Manufacturing"
$ws.Cells.Item(2, 8).Value = 15.39600939973092
$ws.Cells.Item(3, 8).Value = 15.9667660077832
$ws.Cells.Item(4, 8).Value = 16.88003794353869
$ws.Cells.Item(5, 8).Value = 17.07946827549553
$ws.Cells.Item(6, 8).Value = 15.5303734048309

$ws.Cells.Item(1, 9).Value = "Leisure and hospitality"
$ws.Cells.Item(2, 9).Value = 10.95233945413795
$ws.Cells.Item(3, 9).Value = 11.33387519217642
$ws.Cells.Item(4, 9).Value = 10.60218113532374
$ws.Cells.Item(5, 9).Value = 10.99741798254191
$ws.Cells.Item(6, 9).Value = 10.75608397296303

$ws.Cells.Item(1, 10).Value = "Information"
$ws.Cells.Item(2, 10).Value = 0.3572362342451978
$ws.Cells.Item(3, 10).Value = 0.30705315638765
$ws.Cells.Item(4, 10).Value = 0.2750118404046781
$ws.Cells.Item(5, 10).Value = 0.2773262380082226
$ws.Cells.Item(6, 10).Value = 0.3808285405047073

$ws.Cells.Item(1, 11).Value = "Financial activities"
$ws.Cells.Item(2, 11).Value = 3.764049835090717
$ws.Cells.Item(3, 11).Value = 1.968752822226655
$ws.Cells.Item(4, 11).Value = 2.029397820178243
$ws.Cells.Item(5, 11).Value = 1.998661182437321
$ws.Cells.Item(6, 11).Value = 2.266394205011207

$ws.Cells.Item(1, 12).Value = "Education and health services"
$ws.Cells.Item(2, 12).Value = 18.34974295180153
$ws.Cells.Item(3, 12).Value = 19.39853698230673
$ws.Cells.Item(4, 12).Value = 20.04741584979819
$ws.Cells.Item(5, 12).Value = 18.98249973508321
$ws.Cells.Item(6, 12).Value = 18.4469626686568

$ws.Cells.Item(1, 13).Value = "Construction"
$ws.Cells.Item(2, 13).Value = 2.605210453882969
$ws.Cells.Item(3, 13).Value = 2.357084800922737
$ws.Cells.Item(4, 13).Value = 2.399241348154648
$ws.Cells.Item(5, 13).Value = 2.72544711267099
$ws.Cells.Item(6, 13).Value = 2.963031768026941

# --- Step 3: update chart series colors (by position, idx 1..12) ---
$co = $ws.ChartObjects()
$chart = $co.Item(1).Chart
$sc = $chart.SeriesCollection()
$sc.Item(1).Format.Fill.ForeColor.RGB = HexToRgb "68AFFC"
$sc.Item(2).Format.Fill.ForeColor.RGB = HexToRgb "4233A6"
$sc.Item(3).Format.Fill.ForeColor.RGB = HexToRgb "85E5DD"
$sc.Item(4).Format.Fill.ForeColor.RGB = HexToRgb "2A6866"
$sc.Item(5).Format.Fill.ForeColor.RGB = HexToRgb "66DE78"
$sc.Item(6).Format.Fill.ForeColor.RGB = HexToRgb "15974D"
$sc.Item(7).Format.Fill.ForeColor.RGB = HexToRgb "B4D170"
$sc.Item(8).Format.Fill.ForeColor.RGB = HexToRgb "683C00"
$sc.Item(9).Format.Fill.ForeColor.RGB = HexToRgb "CA7E54"
$sc.Item(10).Format.Fill.ForeColor.RGB = HexToRgb "821F48"
$sc.Item(11).Format.Fill.ForeColor.RGB = HexToRgb "F65B68"
$sc.Item(12).Format.Fill.ForeColor.RGB = HexToRgb "EBCECB"

# --- Step 4: remove major gridlines on value axis ---
$valAx = $chart.Axes(2)
$valAx.HasMajorGridlines = $False

# --- Step 5: reposition/resize the chart ---
$obj = $co.Item(1)
$fromLeft = $ws.Cells.Item(8, 1).Left
$fromTop = $ws.Cells.Item(8, 1).Top
$toLeft = $ws.Cells.Item(36, 5).Left + 75
$toTop = $ws.Cells.Item(36, 5).Top + 12
$obj.Left = $fromLeft
$obj.Top = $fromTop
$obj.Width = $toLeft - $fromLeft
$obj.Height = $toTop - $fromTop